$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText "D2" '28.342.87'
Set-CellText "E2" '  -5.37%  '
Set-CellText "D3" '1.840.85'
Set-CellText "E3" '  -5.18%  '
Set-CellText "D4" '1.002'
Set-CellText "E4" '  -0.41%  '
Set-CellText "D5" '330.81'
Set-CellText "E5" '  -1.43%  '
Set-CellText "D6" '1.001'
Set-CellText "E6" '  -0.49%  '
Set-CellText "D7" '0.4605'
Set-CellText "E7" '  -4.91%  '
Set-CellText "E8" '  -6.56%  '
Set-CellText "E9" '  -2.98%  '
Set-CellText "D10" '0.07868'
Set-CellText "E10" '  -3.85%  '
Set-CellText "D11" '0.9679'
Set-CellText "E11" '  -4.83%  '
Set-CellText "E12" '  -7.45%  '
Set-CellText "D13" '1.849.26'
Set-CellText "E13" '  -5.15%  '
Set-CellText "E14" '  -6.21%  '
Set-CellText "D15" '6.933'
Set-CellText "E15" '  -5.17%  '
Set-CellText "D16" '0.06878'
Set-CellText "E16" '  +0.29%  '
Set-CellText "D17" '1.002'
Set-CellText "E17" '  -0.60%  '
Set-CellText "D18" '86.95'
Set-CellText "E18" '  -4.71%  '
Set-CellText "D19" '0.000009972'
Set-CellText "E19" '  -3.91%  '
Set-CellText "D20" '16.93'
Set-CellText "E20" '  -5.07%  '
Set-CellText "D21" '1.001'
Set-CellText "E21" '  -0.54%  '
Set-CellText "D22" '28.377.67'
Set-CellText "E22" '  -5.28%  '
Set-CellText "E23" '  -5.53%  '
Set-CellText "D24" '11.00'
Set-CellText "E24" '  -7.66%  '
Set-CellText "D25" '2.154'
Set-CellText "D26" '2.068.06'
Set-CellText "E26" '  -4.82%  '
Set-CellText "D27" '153.82'
Set-CellText "E27" '  -1.86%  '
Set-CellText "D28" '19.22'
Set-CellText "E28" '  -4.38%  '
Set-CellText "D29" '5.784'
Set-CellText "E29" '  -13.17%  '
Set-CellText "D30" '1.987'
Set-CellText "E30" '  -5.73%  '
Set-CellText "D31" '116.90'
Set-CellText "E31" '  -3.67%  '
Set-CellText "D32" '0.9427'
Set-CellText "E32" '  -6.88%  '
Set-CellText "D33" '0.09306'
Set-CellText "E33" '  -3.48%  '
Set-CellText "D34" '5.290'
Set-CellText "E34" '  -5.56%  '
Set-CellText "D35" '3.444'
Set-CellText "E35" '  -2.92%  '
Set-CellText "E36" '  -6.97%  '
Set-CellText "D37" '0.06032'
Set-CellText "E37" '  -8.40%  '
Set-CellText "D38" '0.02157'
Set-CellText "E38" '  -5.85%  '
Set-CellText "D39" '1.154'
Set-CellText "E39" '  -4.91%  '
Set-CellText "D40" '1.000'
Set-CellText "E40" '  -0.53%  '
Set-CellText "D41" '7.608'
Set-CellText "E41" '  -5.06%  '
Set-CellText "D42" '0.5637'
Set-CellText "E42" '  -5.75%  '
Set-CellText "D43" '10.02'
Set-CellText "E43" '  -6.80%  '
Set-CellText "D44" '0.1786'
Set-CellText "E44" '  -3.66%  '
Set-CellText "D45" '1.236'
Set-CellText "E45" '  -2.78%  '
Set-CellText "D46" '2.278'
Set-CellText "E46" '  -8.68%  '
Set-CellText "E47" '  -5.29%  '
Set-CellText "D48" '0.5302'
Set-CellText "E48" '  -5.05%  '
Set-CellText "D49" '0.07039'
Set-CellText "E49" '  -6.26%  '
Set-CellText "D50" '1.843'
Set-CellText "E50" '  -7.32%  '
Set-CellText "D51" '113.09'
Set-CellText "E51" '  -3.91%  '
